$wb = $excel.ActiveWorkbook

# Sheets "展览" (index 1) and "全部类型" (index 4) both hold the identical
# 17-row (A1:I17) event table and both receive the identical edit:
#   - F2 (想去人数 for row 1) bumps from 1062 to 1066
#   - A brand new event row is inserted as row 3 (pushing every following
#     row down by one); the new row keeps the next sequential id (2)
#   - A handful of the now-shifted rows get their F (想去人数) value bumped,
#     and one of them (景德镇→抚州 "第七届FZ动漫文化节" row) also gets a new
#     cover image URL
#   - A brand new last row (id 17) is appended at the end, carrying what
#     used to be the final row's content unchanged

$targetSheetNames = @("展览", "全部类型")

foreach ($sheetName in $targetSheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # --- update row 2 (first data row) interest count ---
    $ws.Cells.Item(2, 6).Value = 1066

    # --- insert a brand-new row at position 3, shifting rows 3-17 down to 4-18 ---
    $ws.Rows.Item(3).Insert()

    # the inserted row's id cell (A3) picks up the header's bold/border
    # style by default; restore the plain numbered-row style used by every
    # other id cell in column A (copy format from the row right below, which
    # kept the correct style when it shifted down from the old row 3)
    $ws.Cells.Item(4, 1).Copy()
    $ws.Cells.Item(3, 1).PasteSpecial(-4122)

    # --- populate the newly inserted row 3 with the new event ---
    # (B3 holds a date-shaped string like the other rows in this column, not
    # a real date value, so force text formatting first to stop the engine
    # from auto-converting it into a date serial number)
    $ws.Cells.Item(3, 2).NumberFormat = "@"
    $ws.Cells.Item(3, 1).Value = 2
    $ws.Cells.Item(3, 2).Value = "2024-03-10"
    $ws.Cells.Item(3, 3).Value = "抚州·宅舞联萌·随舞动漫派对（免费活动)"
    $ws.Cells.Item(3, 4).Value = "复兴大道333号 华章保利拾光荟购物广场"
    $ws.Cells.Item(3, 5).Value = "2024.03.10 14:00-03.10 18:00"
    $ws.Cells.Item(3, 6).Value = 2
    $ws.Cells.Item(3, 7).Value = 22.33
    $ws.Cells.Item(3, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82373"
    $ws.Cells.Item(3, 9).Value = "//i1.hdslb.com/bfs/openplatform/202403/UfqnH2p81709691112608.jpeg"

    # --- tweak the "想去人数" (F) values of several rows that shifted down ---
    $ws.Cells.Item(5, 6).Value = 3058    # 江西·ShiningStaR动漫游戏文化节5th: 3055 -> 3058
    $ws.Cells.Item(7, 6).Value = 2235    # 南昌·AP动漫游戏嘉年华: 2218 -> 2235
    $ws.Cells.Item(9, 6).Value = 113     # 南昌·AP动漫游戏 嘉年华内场票: 110 -> 113
    $ws.Cells.Item(10, 6).Value = 1054   # 南昌·CM01动漫游戏博览会: 1048 -> 1054
    $ws.Cells.Item(12, 6).Value = 39     # 新余·文旅国漫嘉年华暨BM次元盛典: 37 -> 39
    $ws.Cells.Item(13, 6).Value = 262    # 赣州·第三届半夏动漫展: 256 -> 262
    $ws.Cells.Item(14, 6).Value = 293    # 赣州·赣次元·归来国风动漫节: 93 -> 293
    $ws.Cells.Item(15, 6).Value = 5      # 抚州·第七届FZ动漫文化节: 2 -> 5
    $ws.Cells.Item(16, 6).Value = 97     # 南昌·原X穹X崩only: 96 -> 97

    # --- the cover image for 赣州·赣次元·归来国风动漫节 (now row 14) also changed ---
    $ws.Cells.Item(14, 9).Value = "//i1.hdslb.com/bfs/openplatform/202403/fIehikk51709705287036.jpeg"

    # --- the old last row (南昌·DSL国风动漫游戏嘉年华) already shifted down to row
    #     18 by the insert above, with all its B-I content intact; it just
    #     needs its sequential id bumped from 16 to 17 ---
    $ws.Cells.Item(18, 1).Value = 17
}
